# ---------------------------------------------------------------------------
# Commit: "feat: add 2022-Q4 data"
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (pushing every
#    existing quarter sheet one slot later) and fill it with the new
#    quarter's fund-holding table.
# 2. Insert a new row into the "总计" (totals) sheet for 2022-Q4
#    (13 holdings, 7.97亿元), shifting the existing quarterly rows down by
#    one and renumbering the running index in column A.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q4" worksheet -----------------------------
# Copy an existing quarter sheet (same column layout / header style / index
# style) so formatting is inherited "for free", drop the copy right after
# "总计", rename it, then overwrite its cell values in place.
$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2020-Q4")
$template.Copy($null, $total)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template ("2020-Q4") has 14 data rows, the new quarter only has 13 -
# drop the trailing row once the values below are written.
$newSheet.Rows.Item(15).Delete()

# Header row (only D1 differs between quarters: "基金金额" -> "基金规模")
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# code, name, scale, position%, ratio, marketValue, rank
$q4Data = @(
    @("001645", "国泰大健康股票A",               "36.92", "94.51", "5.86", "2.1635", 10),
    @("009274", "融通健康产业灵活配置混合C",       "49.79", "94.21", "2.66", "1.3244", 8),
    @("020001", "国泰金鹰增长灵活配置混合",         "20.52", "94.14", "5.62", "1.1532", 10),
    @("000727", "融通健康产业灵活配置混合A",       "41.50", "94.21", "2.66", "1.1039", 8),
    @("009805", "国泰医药健康股票A",               "10.35", "93.35", "8.97", "0.9284", 5),
    @("011321", "国泰大健康股票C",                 "5.56",  "94.51", "5.86", "0.3258", 10),
    @("011335", "银河医药健康混合A",               "8.18",  "92.65", "3.83", "0.3133", 10),
    @("160215", "国泰价值经典灵活配置混合（LOF）", "5.95",  "94.08", "5.09", "0.3029", 10),
    @("008370", "国泰研究精选两年持有期混合",       "2.79",  "93.97", "6.31", "0.1760", 8),
    @("519673", "银河康乐股票A",                   "2.15",  "93.79", "4.58", "0.0985", 7),
    @("011326", "国泰医药健康股票C",               "0.79",  "93.35", "8.97", "0.0709", 5),
    @("016018", "银河康乐股票C",                   "0.10",  "93.79", "4.58", "0.0046", 7),
    @("015666", "银河医药健康混合C",               "0.02",  "92.65", "3.83", "0.0008", 10)
)

$newSheet.Range("B2:G14").NumberFormat = "@"

$r = 2
foreach ($row in $q4Data) {
    $newSheet.Cells.Item($r,1).Value = $r - 2
    $newSheet.Cells.Item($r,2).Value = $row[0]
    $newSheet.Cells.Item($r,3).Value = $row[1]
    $newSheet.Cells.Item($r,4).Value = $row[2]
    $newSheet.Cells.Item($r,5).Value = $row[3]
    $newSheet.Cells.Item($r,6).Value = $row[4]
    $newSheet.Cells.Item($r,7).Value = $row[5]
    $newSheet.Cells.Item($r,8).Value = $row[6]
    $r = $r + 1
}

# --- Step 2: insert the 2022-Q4 row into the "总计" (totals) sheet ----------
$total.Rows.Item(2).Insert()
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 13
$total.Cells.Item(2,4).Value = 7.97

# Renumber the running index in column A for every row below the new one.
for ($row = 3; $row -le 10; $row++) {
    $total.Cells.Item($row,1).Value = $row - 2
}
